# Actualización automática del tracker
# Appends 3 new result rows (38-40) to the results tracker sheet, matching
# the existing column layout:
# event_id | fecha | jugador_A | jugador_B | pronostico | cuota | resultado | profit
# These are newly-registered picks, so "resultado"/"profit" (G/H) are still
# blank, same as every other pending row already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Cells.Item(38, 1).Value = 14343570
$ws.Cells.Item(38, 2).Value = "'2025-08-04"
$ws.Cells.Item(38, 3).Value = "Toby Alex Kodat"
$ws.Cells.Item(38, 4).Value = "Rafael Jodar"
$ws.Cells.Item(38, 5).Value = "Gana Toby Alex Kodat"
$ws.Cells.Item(38, 6).Value = 3

# Row 39
$ws.Cells.Item(39, 1).Value = 14344481
$ws.Cells.Item(39, 2).Value = "'2025-08-04"
$ws.Cells.Item(39, 3).Value = "Fabrizio Andaloro"
$ws.Cells.Item(39, 4).Value = "Maximus Jones"
$ws.Cells.Item(39, 5).Value = "Gana Maximus Jones"
$ws.Cells.Item(39, 6).Value = 2.75

# Row 40
$ws.Cells.Item(40, 1).Value = 14344408
$ws.Cells.Item(40, 2).Value = "'2025-08-04"
$ws.Cells.Item(40, 3).Value = "Mariano Kestelboim"
$ws.Cells.Item(40, 4).Value = "Michael Vrbensky"
$ws.Cells.Item(40, 5).Value = "Gana Michael Vrbensky"
$ws.Cells.Item(40, 6).Value = 1.83

# Keep "resultado"/"profit" blank (same shape as the rest of the sheet's
# pending rows) by replicating the existing blank G/H cell from row 37.
$ws.Range("G37:H37").Copy($ws.Range("G38:H40"))
